# The deck ships two theme parts:
#   ppt/theme/theme1.xml  - used by the slide master (currently "Integral")
#   ppt/theme/theme2.xml  - used by the notes master  (currently "Office Theme")
# The target edit swaps the two themes' color schemes so that theme1.xml
# ends up with the Office Theme palette (and theme2.xml ends up with the
# Integral palette).
#
# PowerPoint's ThemeColorScheme exposes the 12 scheme colors in the
# standard ColorScheme index order:
#   1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3 8=accent4
#   9=accent5 10=accent6 11=hlink 12=folHlink
#
# ColorFormat.RGB takes/returns a VBA-style BGR-packed integer
# (0x00BBGGRR), i.e. the reverse byte order of the hex RRGGBB values
# found in the OOXML <a:srgbClr val="RRGGBB"/> markup.

$p = $ppt.ActivePresentation

# "Office Theme" palette (RRGGBB -> BGR int) to apply to the slide master's
# theme (theme1.xml).
$officeColors = @(
    0x000000,  # dk1      000000
    0xFFFFFF,  # lt1      FFFFFF
    0x6A5444,  # dk2      44546A
    0xE6E6E7,  # lt2      E7E6E6
    0xD59B5B,  # accent1  5B9BD5
    0x317DED,  # accent2  ED7D31
    0xA5A5A5,  # accent3  A5A5A5
    0x00C0FF,  # accent4  FFC000
    0xC47244,  # accent5  4472C4
    0x47AD70,  # accent6  70AD47
    0xC16305,  # hlink    0563C1
    0x724F95   # folHlink 954F72
)

$masterScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $masterScheme.Item($i).RGB = $officeColors[$i - 1]
}
